# Edit: generate named schema:Place's
#
# Every other CmsLocation row (the "odd"-indexed work/freestanding-work
# locations) is pulled out into a new, dedicated SchemaPlace sheet that
# holds just the @id of a named place. The remaining ("even"-indexed)
# locations stay behind as plain (unnamed) CmsLocation rows.
#
# The new SchemaPlace sheet is inserted right after CmsWork / right before
# SchemaCreativeWork, matching the workbook's existing sheet ordering
# convention of Cms<Foo> followed by Schema<Foo>.

$wb = $excel.ActiveWorkbook

# --- Trim CmsLocation down to the locations that stay unnamed ---
$loc = $wb.Worksheets.Item("CmsLocation")

$loc.Range("A3").Value = "http://example.com/collection0/work2Location"
$loc.Range("A4").Value = "http://example.com/collection1/work4Location"
$loc.Range("A5").Value = "http://example.com/collection1/work6Location"
$loc.Range("A6").Value = "http://example.com/freestandingwork8Location"
$loc.Range("A7").Value = "http://example.com/freestandingwork10Location"

# Rows 8-13 held the locations that moved to the new SchemaPlace sheet
$loc.Rows("8:13").Delete() | Out-Null

# --- Create the SchemaPlace sheet, positioned before SchemaCreativeWork ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "SchemaPlace"

$target = $wb.Worksheets.Item("SchemaCreativeWork")
$newSheet.Move($target)

# Re-fetch by name: the Move() call invalidates the old object handle
$place = $wb.Worksheets.Item("SchemaPlace")

$place.Range("A1").Value = "@id"
$place.Range("A2").Value = "http://example.com/collection0/work1Location"
$place.Range("A3").Value = "http://example.com/collection0/work3Location"
$place.Range("A4").Value = "http://example.com/collection1/work5Location"
$place.Range("A5").Value = "http://example.com/collection1/work7Location"
$place.Range("A6").Value = "http://example.com/freestandingwork9Location"
$place.Range("A7").Value = "http://example.com/freestandingwork11Location"
